{"js": "/*\n * Update the worksheet date and the 25 division-answer cells.\n * Source document has exactly 26 non-empty paragraphs in document order:\n *   paragraph 0            -> the centered date line\n *   paragraphs 1..25       -> the answer cells of the 5x5 table (row major,\n *                             skipping the blank spacer rows)\n * We walk body.paragraphs in order, skip the empty ones (blank table rows),\n * and overwrite the text of each of the 26 remaining paragraphs with the\n * corresponding new value - this is robust even though a couple of the new\n * values collide with each other / with older values elsewhere in the\n * document.\n */\nconst expectedOld = [\"2023-09-25 Monday\", \"27\u00f78=3, 3\", \"49\u00f79=5, 4\", \"78\u00f76=13, 0\", \"52\u00f73=17, 1\", \"97\u00f77=13, 6\", \"21\u00f77=3, 0\", \"99\u00f76=16, 3\", \"38\u00f74=9, 2\", \"98\u00f77=14, 0\", \"24\u00f79=2, 6\", \"31\u00f77=4, 3\", \"50\u00f78=6, 2\", \"97\u00f75=19, 2\", \"66\u00f72=33, 0\", \"51\u00f75=10, 1\", \"45\u00f72=22, 1\", \"80\u00f75=16, 0\", \"96\u00f74=24, 0\", \"32\u00f74=8, 0\", \"21\u00f73=7, 0\", \"34\u00f74=8, 2\", \"24\u00f73=8, 0\", \"29\u00f78=3, 5\", \"84\u00f79=9, 3\", \"99\u00f77=14, 1\"];\nconst newValues = [\"2023-09-26 Tuesday\", \"41\u00f76=6, 5\", \"12\u00f78=1, 4\", \"39\u00f74=9, 3\", \"12\u00f72=6, 0\", \"70\u00f77=10, 0\", \"54\u00f79=6, 0\", \"20\u00f72=10, 0\", \"80\u00f79=8, 8\", \"79\u00f73=26, 1\", \"99\u00f79=11, 0\", \"73\u00f72=36, 1\", \"94\u00f77=13, 3\", \"72\u00f78=9, 0\", \"52\u00f75=10, 2\", \"83\u00f78=10, 3\", \"79\u00f73=26, 1\", \"80\u00f76=13, 2\", \"53\u00f75=10, 3\", \"38\u00f74=9, 2\", \"21\u00f74=5, 1\", \"12\u00f76=2, 0\", \"52\u00f78=6, 4\", \"85\u00f77=12, 1\", \"79\u00f72=39, 1\", \"31\u00f75=6, 1\"];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet cursor = 0;\nfor (const paragraph of paragraphs.items) {\n  if (cursor >= newValues.length) break;\n  const current = paragraph.text;\n  if (current === \"\") continue; // blank spacer-row paragraph, leave untouched\n\n  if (current !== expectedOld[cursor]) {\n    throw new Error(\n      `Unexpected paragraph text at position ${cursor}: found ${JSON.stringify(current)}, ` +\n      `expected ${JSON.stringify(expectedOld[cursor])}`\n    );\n  }\n\n  paragraph.insertText(newValues[cursor], Word.InsertLocation.replace);\n  cursor++;\n}\n\nawait context.sync();\n\nif (cursor !== newValues.length) {\n  throw new Error(`Only updated ${cursor} of ${newValues.length} expected paragraphs`);\n}\n", "ps1": "# Update the worksheet date and the 25 division-answer cells.\n# Each (old, new) pair below corresponds 1:1, in document order, to one\n# of the 26 non-empty paragraphs (the date line plus the 25 filled table\n# cells). We run Find/Replace (first-match only, case sensitive, no\n# wildcards) once per pair; since every \"old\" value is unique in the\n# document at the moment it is searched for, each call touches exactly\n# the intended run.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2023-09-25 Monday\", \"2023-09-26 Tuesday\"),\n    @(\"27\u00f78=3, 3\", \"41\u00f76=6, 5\"),\n    @(\"49\u00f79=5, 4\", \"12\u00f78=1, 4\"),\n    @(\"78\u00f76=13, 0\", \"39\u00f74=9, 3\"),\n    @(\"52\u00f73=17, 1\", \"12\u00f72=6, 0\"),\n    @(\"97\u00f77=13, 6\", \"70\u00f77=10, 0\"),\n    @(\"21\u00f77=3, 0\", \"54\u00f79=6, 0\"),\n    @(\"99\u00f76=16, 3\", \"20\u00f72=10, 0\"),\n    @(\"38\u00f74=9, 2\", \"80\u00f79=8, 8\"),\n    @(\"98\u00f77=14, 0\", \"79\u00f73=26, 1\"),\n    @(\"24\u00f79=2, 6\", \"99\u00f79=11, 0\"),\n    @(\"31\u00f77=4, 3\", \"73\u00f72=36, 1\"),\n    @(\"50\u00f78=6, 2\", \"94\u00f77=13, 3\"),\n    @(\"97\u00f75=19, 2\", \"72\u00f78=9, 0\"),\n    @(\"66\u00f72=33, 0\", \"52\u00f75=10, 2\"),\n    @(\"51\u00f75=10, 1\", \"83\u00f78=10, 3\"),\n    @(\"45\u00f72=22, 1\", \"79\u00f73=26, 1\"),\n    @(\"80\u00f75=16, 0\", \"80\u00f76=13, 2\"),\n    @(\"96\u00f74=24, 0\", \"53\u00f75=10, 3\"),\n    @(\"32\u00f74=8, 0\", \"38\u00f74=9, 2\"),\n    @(\"21\u00f73=7, 0\", \"21\u00f74=5, 1\"),\n    @(\"34\u00f74=8, 2\", \"12\u00f76=2, 0\"),\n    @(\"24\u00f73=8, 0\", \"52\u00f78=6, 4\"),\n    @(\"29\u00f78=3, 5\", \"85\u00f77=12, 1\"),\n    @(\"84\u00f79=9, 3\", \"79\u00f72=39, 1\"),\n    @(\"99\u00f77=14, 1\", \"31\u00f75=6, 1\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $found = $range.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Find/Replace failed: could not locate text [$oldText]\"\n    }\n}\n\n"}
